$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Message to the tutor"
$ws.Range("A2").Value = "Deji"
$ws.Range("B2").Value = "I do love your teaching!!!"

$ws.Range("A3").Select()
